$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.8634222919937206
$ws.Range("C2").Value = 0.9458297506448839
$ws.Range("D2").Value = 0.9027492819039803
$ws.Range("E2").Value = 1163

$ws.Range("B3").Value = 0.9581320450885669
$ws.Range("C3").Value = 0.9239130434782609
$ws.Range("D3").Value = 0.9407114624505928
$ws.Range("E3").Value = 644

$ws.Range("B4").Value = 0.8534031413612565
$ws.Range("C4").Value = 0.8402061855670103
$ws.Range("D4").Value = 0.8467532467532468
$ws.Range("E4").Value = 776

$ws.Range("B5").Value = 0.8256227758007118
$ws.Range("C5").Value = 0.6498599439775911
$ws.Range("D5").Value = 0.7272727272727274
$ws.Range("E5").Value = 357

$ws.Range("B6").Value = 0.8772108843537415
$ws.Range("C6").Value = 0.8772108843537415
$ws.Range("D6").Value = 0.8772108843537415
$ws.Range("E6").Value = 0.8772108843537415

$ws.Range("B7").Value = 0.8751450635610639
$ws.Range("C7").Value = 0.8399522309169365
$ws.Range("D7").Value = 0.8543716795951368

$ws.Range("B8").Value = 0.8769337861506541
$ws.Range("C8").Value = 0.8772108843537415
$ws.Range("D8").Value = 0.8749770339419707
